$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2..102 (Cluster name, Active cases, Exist)
$data = @(
    @("ABD Group 512 Melbourne Road Construction Site Spotswood", 5, "old"),
    @("Acquire BPO Southbank", 6, "new"),
    @("Acquire BPO Southbank", 7, "old"),
    @("Al Haj Halal Meats Glenroy", 31, "new"),
    @("Al Haj Halal Meats Glenroy", 34, "old"),
    @("Al-Taqwa College Truganina", 5, "old"),
    @("Al-Taqwa College Truganina", 9, "new"),
    @("Amiga Montessori Craigieburn", 28, "new"),
    @("Anglicare Home Princes Hwy Dandenong", 6, "new"),
    @("Campbellfield Ford Complex Vaccination Clinic Campbellfield", 12, "new"),
    @("Caroline Springs Police Station", 11, "new"),
    @("Caroline Springs Police Station", 12, "old"),
    @("Cedars Medical Clinic Coburg", 28, "old"),
    @("Cedars Medical Clinic Coburg", 39, "new"),
    @("Cedars Medical Clinic Coburg", 41, "old"),
    @("Chemist Warehouse Campbellfield DC", 5, "old"),
    @("Chemist Warehouse Campbellfield DC", 6, "new"),
    @("Chemist Warehouse Fillo Drive Somerton", 5, "old"),
    @("Chemist Warehouse Fillo Drive Somerton", 7, "new"),
    @("City of Hobsons Bay Community", 5, "new"),
    @("City of Moreland Community", 5, "new"),
    @("City of Moreland Community", 6, "old"),
    @("Classy Cabinets and Kitchens Craigieburn", 9, "new"),
    @("Classy Cabinets and Kitchens Craigieburn", 10, "old"),
    @("Coles Broadmeadows Central Shopping Centre", 8, "new"),
    @("Coles Broadmeadows Central Shopping Centre", 9, "old"),
    @("Coles Campbellfield Plaza Campbellfield", 9, "old"),
    @("Coles Campbellfield Plaza Campbellfield", 12, "new"),
    @("Coles Coburg North Village", 26, "old"),
    @("Coles Coburg North Village", 29, "new"),
    @("Coles Pakenham Place Shopping", 11, "new"),
    @("Coles Pakenham Place Shopping Centre", 11, "old"),
    @("Coles Roxburgh Village Roxburgh Park", 10, "new"),
    @("Coles Roxburgh Village Roxburgh Park", 11, "old"),
    @("Community Kids Meadow Heights", 11, "new"),
    @("Community Kids Meadow Heights", 12, "old"),
    @("Costco Wholesale Epping", 24, "old"),
    @("Costco Wholesale Epping", 28, "new"),
    @("Crusader Caravans Epping", 14, "old"),
    @("Crusader Caravans Epping", 21, "new"),
    @("Direct Freight Express Cambellfield", 13, "old"),
    @("Direct Freight Express Campbellfield", 13, "new"),
    @("Don Watson Coldstore Derrimut", 5, "new"),
    @("Epworth Healthcare Epworth Richmond Emergency Department", 5, "old"),
    @("Epworth Healthcare Epworth Richmond Emergency Department", 6, "new"),
    @("Fine Food Holdings Pty Ltd Dandenong South", 6, "new"),
    @("General foods Campbellfield", 12, "new"),
    @("Goodstart Early Learning Altona", 5, "old"),
    @("Goodstart Early Learning Altona", 7, "new"),
    @("Green Leaves Early Learning Centre Highlands Craigieburn", 7, "old"),
    @("Green Leaves Early Learning Centre Highlands Craigieburn", 8, "new"),
    @("Hamilton Marino 236 Jasper Road McKinnon", 12, "old"),
    @("Hamilton Marino 236 Jasper Road McKinnon", 13, "new"),
    @("Health Care Providers Association South Melbourne", 9, "new"),
    @("Health Care Providers Association South Melbourne", 13, "old"),
    @("Islamic College of Melbourne Tarneit", 5, "old"),
    @("KFC Fawkner", 7, "new"),
    @("KFC Fawkner", 8, "old"),
    @("Kids House Early Learning Cheltenham", 11, "new"),
    @("Kids House Early Learning Cheltenham", 12, "old"),
    @("Learning Nest Early Learning Centre Meadow Heights", 5, "old"),
    @("McDonalds Thomastown II", 5, "new"),
    @("McDonalds Thomastown II", 7, "old"),
    @("Melton Police Station Melton", 5, "old"),
    @("Montessori Beginnings Greenvale", 5, "old"),
    @("MyCentre Childcare Broadmeadows", 16, "new"),
    @("MyCentre Childcare Broadmeadows", 17, "old"),
    @("Newbury Child and Community Centre Craigieburn", 5, "old"),
    @("Nido Early School Moonee Ponds", 14, "old"),
    @("Nido Early School Moonee Ponds", 15, "new"),
    @("Northern Health Northern Hospital Epping Emergency Department Tier 1B", 47, "old"),
    @("Northern Health Northern Hospital Epping Emergency Department Tier 1B", 48, "new"),
    @("Northern Health The Northern Hospital Epping", 5, "old"),
    @("Northern Health The Northern Hospital Epping", 8, "new"),
    @("Oporto Coolaroo", 12, "new"),
    @("Oporto Coolaroo", 13, "old"),
    @("Our Lady Help of Christian's Primary School Brunswick East", 8, "old"),
    @("Our Lady Help of Christian's Primary School Brunswick East", 10, "new"),
    @("Paisley Park Early Learning Centre Bundoora", 7, "new"),
    @("Paisley Park Early Learning Centre Bundoora", 8, "old"),
    @("Panorama Construction Site Whitehorse Rd Box Hill", 13, "new"),
    @("Panorama Construction Site Whitehorse Rd Box Hill", 17, "old"),
    @("Richmond Quarter 261-271 Bridge Road Construction Site", 11, "new"),
    @("Richmond Quarter 261-271 Bridge Road Construction Site Richmond", 10, "old"),
    @("Sacca's Fruit World Broadmeadows Central Shopping Centre", 5, "new"),
    @("St Margaret's Primary School Maribyrnong", 11, "new"),
    @("Tek Foods Somerton", 12, "old"),
    @("Tek Foods Somerton", 14, "new"),
    @("The Homestead Child and Family Centre Roxburgh Park", 11, "new"),
    @("The Homestead Child and Family Centre Roxburgh Park", 12, "old"),
    @("The Huntly-Goornong Rail Works", 5, "new"),
    @("The Royal Children's Hospital Melbourne Emergency Department Parkville Tier 1B", 11, "old"),
    @("The Royal Children's Hospital Melbourne Emergency Department Parkville Tier 1B", 17, "new"),
    @("Unilodge College Square Student Accommodation 570 Lygon Street", 14, "new"),
    @("Unilodge College Square Student Accommodation 570 Lygon Street Carlton", 14, "old"),
    @("Werribee Mercy Hospital Emergency Department", 8, "old"),
    @("Werribee Mercy Hospital Emergency Department", 10, "new"),
    @("Western Health Sunshine Hospital Emergency Department", 8, "old"),
    @("Western Health Sunshine Hospital Emergency Department", 11, "new"),
    @("Yara Childcare Centre Truganina", 9, "new"),
    @("Yara Childcare Centre Truganina", 10, "old")
)

$rowCount = $data.Count
$arr = New-Object 'object[,]' $rowCount,3
for ($i = 0; $i -lt $rowCount; $i++) {
    $arr[$i, 0] = $data[$i][0]
    $arr[$i, 1] = $data[$i][1]
    $arr[$i, 2] = $data[$i][2]
}

$startRow = 2
$endRow = $startRow + $rowCount - 1
$rng = $ws.Range("A${startRow}:C${endRow}")
$rng.Value = $arr
